# Natmi following Dr Hou advice
# Update Ligand-expressing cells / Receptor-expressing cells counts (1 -> 3)
# and the downstream derived expression/specificity metrics for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=1.659925666666667;  H=4.979777;            I=0.2642801384529059; J=0.2642801384529059; K=3; M=7.446932333333334;  N=22.340797;  O=0.1309571246780483; P=0.1309571246780484; Q=12.36135411802989;  R=111.252187062269;   S=0.03460936704130908; T=0.03460936704130909 }
    3  = @{ E=3; G=1.659925666666667;  H=4.979777;            I=0.2642801384529059; J=0.2642801384529059; K=3; M=35.70141366666667;  N=107.104241; O=0.6278228767838827; P=0.6278228767838828; Q=59.26169288158412;  R=533.3552359342571;  S=0.1659211168003462;  T=0.1659211168003462 }
    4  = @{ E=3; G=1.659925666666667;  H=4.979777;            I=0.2642801384529059; J=0.2642801384529059; K=3; M=13.71707733333333;  N=41.151232;  O=0.2412199985380688; P=0.2412199985380689; Q=22.76932873725156;  R=204.923958635264;   S=0.06374965461125058; T=0.06374965461125061 }
    5  = @{ E=3; G=1.475152;           H=4.425456000000001;  I=0.2348619475123571; J=0.2348619475123571; K=3; M=7.446932333333334;  N=22.340797;  O=0.1309571246780483; P=0.1309571246780484; Q=10.98535712538134;  R=98.86821412843202;  S=0.030756845342505;   T=0.030756845342505 }
    6  = @{ E=3; G=1.475152;           H=4.425456000000001;  I=0.2348619475123571; J=0.2348619475123571; K=3; M=35.70141366666667;  N=107.104241; O=0.6278228767838827; P=0.6278228767838828; Q=52.66501177321068;  R=473.985105958896;   S=0.1474517035342733;  T=0.1474517035342733 }
    7  = @{ E=3; G=1.475152;           H=4.425456000000001;  I=0.2348619475123571; J=0.2348619475123571; K=3; M=13.71707733333333;  N=41.151232;  O=0.2412199985380688; P=0.2412199985380689; Q=20.23477406242134;  R=182.112966561792;   S=0.05665339863557878; T=0.05665339863557879 }
    8  = @{ E=3; G=3.145854666666667;  H=9.437564;            I=0.5008579140347369; J=0.500857914034737;  K=3; M=7.446932333333334;  N=22.340797;  O=0.1309571246780483; P=0.1309571246780484; Q=23.42696683316756;  R=210.842701498508;   S=0.06559091229423426; T=0.06559091229423429 }
    9  = @{ E=3; G=3.145854666666667;  H=9.437564;            I=0.5008579140347369; J=0.500857914034737;  K=3; M=35.70141366666667;  N=107.104241; O=0.6278228767838827; P=0.6278228767838828; Q=112.3114587898804;  R=1010.803129108924;  S=0.3144500564492632;  T=0.3144500564492633 }
    10 = @{ E=3; G=3.145854666666667;  H=9.437564;            I=0.5008579140347369; J=0.500857914034737;  K=3; M=13.71707733333333;  N=41.151232;  O=0.2412199985380688; P=0.2412199985380689; Q=43.15193174209423;  R=388.367385678848;   S=0.1208169452912394;  T=0.1208169452912395 }
}

foreach ($rowNum in $data.Keys) {
    $row = $data[$rowNum]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$rowNum").Value = $row[$col]
    }
}
